# Updates average_county_temperature (column I), worst_ashp_cop (column N)
# and best_ashp_cop (column O) for a subset of rows, reflecting refreshed
# NOAA county-temperature data merged back into the longform dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> hashtable of column letter -> new numeric value
$updates = @{
    4  = @{ I = 12.51681286549706;  N = 1.782371783972741; O = 1.939565227172176 }
    5  = @{ I = 15.74228395061728;  N = 1.837513876759573; O = 2.005936573945218 }
    6  = @{ I = 1.925925925925943;  N = 1.62249843161857;  O = 1.749494516792324 }
    7  = @{ I = 1.925925925925943 }
    8  = @{ I = 12.66820987654322;  N = 1.784885911058073; O = 1.942582169301264 }
    9  = @{ I = 12.66820987654322 }
    11 = @{ I = -3.222222222222223; N = 1.554711451758341; O = 1.669946025515211 }
    15 = @{ I = 1.925925925925943;  N = 1.62249843161857;  O = 1.749494516792324 }
    16 = @{ I = 20.68981481481483;  N = 1.929056920423291; O = 2.117059768804106 }
    17 = @{ I = 20.68981481481483 }
    18 = @{ I = 14.96875;           N = 1.8239809580482;   O = 1.989608681354817 }
    19 = @{ I = 1.925925925925943;  N = 1.62249843161857;  O = 1.749494516792324 }
    20 = @{ I = -3.222222222222223; N = 1.554711451758341; O = 1.669946025515211 }
    21 = @{ I = -3.222222222222223 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
